$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(807, 1).Value = "Atlanta Hawks"
$ws.Cells.Item(807, 2).Value = 99
$ws.Cells.Item(807, 3).Value = "Charlotte Hornets"
$ws.Cells.Item(807, 4).Value = 122
$ws.Cells.Item(807, 6).Value = 17832
$ws.Cells.Item(807, 7).Value = "Spectrum Center"
$ws.Cells.Item(807, 8).Value = "Charlotte Hornets"
$ws.Cells.Item(807, 9).Value = "Atlanta Hawks"

$ws.Cells.Item(808, 1).Value = "New York Knicks"
$ws.Cells.Item(808, 2).Value = 100
$ws.Cells.Item(808, 3).Value = "Orlando Magic"
$ws.Cells.Item(808, 4).Value = 118
$ws.Cells.Item(808, 6).Value = 17832
$ws.Cells.Item(808, 7).Value = "Amway Center"
$ws.Cells.Item(808, 8).Value = "Orlando Magic"
$ws.Cells.Item(808, 9).Value = "New York Knicks"

$ws.Cells.Item(809, 1).Value = "Miami Heat"
$ws.Cells.Item(809, 2).Value = 109
$ws.Cells.Item(809, 3).Value = "Philadelphia 76ers"
$ws.Cells.Item(809, 4).Value = 104
$ws.Cells.Item(809, 6).Value = 17832
$ws.Cells.Item(809, 7).Value = "Wells Fargo Center"
$ws.Cells.Item(809, 8).Value = "Miami Heat"
$ws.Cells.Item(809, 9).Value = "Philadelphia 76ers"

$ws.Cells.Item(810, 1).Value = "Brooklyn Nets"
$ws.Cells.Item(810, 2).Value = 86
$ws.Cells.Item(810, 3).Value = "Boston Celtics"
$ws.Cells.Item(810, 4).Value = 136
$ws.Cells.Item(810, 6).Value = 17832
$ws.Cells.Item(810, 7).Value = "TD Garden"
$ws.Cells.Item(810, 8).Value = "Boston Celtics"
$ws.Cells.Item(810, 9).Value = "Brooklyn Nets"

$ws.Cells.Item(811, 1).Value = "Chicago Bulls"
$ws.Cells.Item(811, 2).Value = 105
$ws.Cells.Item(811, 3).Value = "Cleveland Cavaliers"
$ws.Cells.Item(811, 4).Value = 108
$ws.Cells.Item(811, 6).Value = 17832
$ws.Cells.Item(811, 7).Value = "Rocket Mortgage Fieldhouse"
$ws.Cells.Item(811, 8).Value = "Cleveland Cavaliers"
$ws.Cells.Item(811, 9).Value = "Chicago Bulls"

$ws.Cells.Item(812, 1).Value = "Indiana Pacers"
$ws.Cells.Item(812, 2).Value = 127
$ws.Cells.Item(812, 3).Value = "Toronto Raptors"
$ws.Cells.Item(812, 4).Value = 125
$ws.Cells.Item(812, 6).Value = 17832
$ws.Cells.Item(812, 7).Value = "Scotiabank Arena"
$ws.Cells.Item(812, 8).Value = "Indiana Pacers"
$ws.Cells.Item(812, 9).Value = "Toronto Raptors"

$ws.Cells.Item(813, 1).Value = "Houston Rockets"
$ws.Cells.Item(813, 2).Value = 113
$ws.Cells.Item(813, 3).Value = "Memphis Grizzlies"
$ws.Cells.Item(813, 4).Value = 121
$ws.Cells.Item(813, 6).Value = 17832
$ws.Cells.Item(813, 7).Value = "FedEx Forum"
$ws.Cells.Item(813, 8).Value = "Memphis Grizzlies"
$ws.Cells.Item(813, 9).Value = "Houston Rockets"

$ws.Cells.Item(814, 1).Value = "Washington Wizards"
$ws.Cells.Item(814, 2).Value = 126
$ws.Cells.Item(814, 3).Value = "New Orleans Pelicans"
$ws.Cells.Item(814, 4).Value = 133
$ws.Cells.Item(814, 6).Value = 17832
$ws.Cells.Item(814, 7).Value = "Smoothie King Center"
$ws.Cells.Item(814, 8).Value = "New Orleans Pelicans"
$ws.Cells.Item(814, 9).Value = "Washington Wizards"

$ws.Cells.Item(815, 1).Value = "San Antonio Spurs"
$ws.Cells.Item(815, 2).Value = 93
$ws.Cells.Item(815, 3).Value = "Dallas Mavericks"
$ws.Cells.Item(815, 4).Value = 116
$ws.Cells.Item(815, 6).Value = 17832
$ws.Cells.Item(815, 7).Value = "American Airlines Center"
$ws.Cells.Item(815, 8).Value = "Dallas Mavericks"
$ws.Cells.Item(815, 9).Value = "San Antonio Spurs"

$ws.Cells.Item(816, 1).Value = "Sacramento Kings"
$ws.Cells.Item(816, 2).Value = 102
$ws.Cells.Item(816, 3).Value = "Denver Nuggets"
$ws.Cells.Item(816, 4).Value = 98
$ws.Cells.Item(816, 6).Value = 17832
$ws.Cells.Item(816, 7).Value = "Ball Arena"
$ws.Cells.Item(816, 8).Value = "Sacramento Kings"
$ws.Cells.Item(816, 9).Value = "Denver Nuggets"

$ws.Cells.Item(817, 1).Value = "Detroit Pistons"
$ws.Cells.Item(817, 2).Value = 100
$ws.Cells.Item(817, 3).Value = "Phoenix Suns"
$ws.Cells.Item(817, 4).Value = 116
$ws.Cells.Item(817, 6).Value = 17832
$ws.Cells.Item(817, 7).Value = "Footprint Center"
$ws.Cells.Item(817, 8).Value = "Phoenix Suns"
$ws.Cells.Item(817, 9).Value = "Detroit Pistons"

$ws.Cells.Item(818, 1).Value = "Los Angeles Lakers"
$ws.Cells.Item(818, 2).Value = 138
$ws.Cells.Item(818, 3).Value = "Utah Jazz"
$ws.Cells.Item(818, 4).Value = 122
$ws.Cells.Item(818, 6).Value = 17832
$ws.Cells.Item(818, 7).Value = "Delta Center"
$ws.Cells.Item(818, 8).Value = "Los Angeles Lakers"
$ws.Cells.Item(818, 9).Value = "Utah Jazz"

$ws.Cells.Item(819, 1).Value = "Los Angeles Clippers"
$ws.Cells.Item(819, 2).Value = 130
$ws.Cells.Item(819, 3).Value = "Golden State Warriors"
$ws.Cells.Item(819, 4).Value = 125
$ws.Cells.Item(819, 6).Value = 17832
$ws.Cells.Item(819, 7).Value = "Chase Center"
$ws.Cells.Item(819, 8).Value = "Los Angeles Clippers"
$ws.Cells.Item(819, 9).Value = "Golden State Warriors"

try {
    $excel.ActiveWindow.ScrollRow = 787
} catch {}
$ws.Range("F808").Select()
